# ADAPT fix costs, cop for TTES and cop for HP (air) based on temperature
#
# Updates the TTES coefficient of performance (column G, header "p_ttes_cop")
# from 5.3 to 12.105 on cell G2 of every yearly sheet (2025, 2030, 2035,
# 2040, 2045, 2050).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("G2").Value = 12.105
}
